$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333, shifting existing rows 333-420 down to 334-421.
$ws.Rows(333).Insert()

# Populate the newly inserted row 333 with the latest week's price record
# (same market/product/category as the surrounding rows, new date + values).
$ws.Range("A333").Value = 3
$ws.Range("B333").Value = "Femacal de La Calera"
$ws.Range("C333").Value = "Coquimbo"
$ws.Range("D333").Value = 44736
$ws.Range("E333").Value = 5
$ws.Range("F333").Value = "Fruta"
$ws.Range("G333").Value = 100108
$ws.Range("H333").Value = "Tropicales y subtropicales"
$ws.Range("I333").Value = 100108002
$ws.Range("J333").Value = "Mango"
$ws.Range("K333").Value = "Sin especificar"
$ws.Range("L333").Value = "Primera"
$ws.Range("M333").Value = 250
$ws.Range("N333").Value = 9000
$ws.Range("O333").Value = 9000
$ws.Range("P333").Value = 9000
$ws.Range("Q333").Value = "$/bandeja 4 kilos"
$ws.Range("R333").Value = "Brasil"
$ws.Range("S333").Value = 2250
$ws.Range("T333").Value = 4
